$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "317.75"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "3.60%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "39.83"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2.42%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.141"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "0.60%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08224"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.55%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "8.345"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "4.57%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9381"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.88%"
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1375"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-6.75%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1984"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.69%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.09118"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "0.17%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03508"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "0.10%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09789"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "0.08%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.001408"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.36%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.006236"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "2.80%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.694"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-2.33%"
$ws.Range("B17").Value = "GateToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.333"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "3.49%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.205"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-6.44%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.3496"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "1.03%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1311"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.56%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.946"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "5.83%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.2449"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.39%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04353"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-0.58%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001228"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-0.68%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004827"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "12.70%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001300"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.12%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0003996"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "-10.15%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02224"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "8.97%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05206"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.62%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007741"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.64%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.009687"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-4.14%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1409"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "4.40%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002049"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-3.48%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.009660"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-2.47%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006616"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "6.96%"
$ws.Range("B48").Value = "CoinbaseStockToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001688"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "5.42%"
$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.002943"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "1.88%"
